$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.07771692647734278, 0.494001664178519, 0.4064399038399985, 0.6375263946222137, 0.6414999648744212, 37)
    3  = @(0.2046287817491858, 0.5768629073926809, 0.7002384789106324, 0.8368025328060572, 0.8229070327852029, 36)
    4  = @(0.09125273495800906, 0.4483188659895491, 0.3552272405021807, 0.5960094298768944, 0.5975810616996375, 35)
    5  = @(0.25374129848085, 0.6003975769414247, 0.8008700048052377, 0.8949134063166323, 0.8710930469124941, 34)
    6  = @(0.1232924819113544, 0.4343337923113413, 0.3678743112519078, 0.606526430794164, 0.6030707379774554, 33)
    7  = @(0.2862287335415685, 0.6098011728780182, 0.8494587501165254, 0.9216608650238576, 0.8901072713674041, 32)
    8  = @(0.09860540898821402, 0.3873856125923802, 0.2559223495544988, 0.5058876847231002, 0.5043867236243681, 31)
    9  = @(0.3607703471504674, 0.6436156715191851, 0.9577333858632414, 0.9786385368782702, 0.9252649997741594, 30)
    10 = @(0.1695398884498732, 0.3837015019462682, 0.2676768073393338, 0.5173749195113093, 0.4974599041038531, 29)
    11 = @(0.3374360160835732, 0.5955736213918028, 0.7831839435432164, 0.8849768039577175, 0.8331329670364563, 28)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
